$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Paragraph 1: "Fitnesse on Node.js" -> "Acceptance Tests running on
# Node with Decaf and Coffeescript"
#
# We rebuild the paragraph (instead of just swapping the run text) so
# that the stray w:proofErr spell-check markers that decorated
# "Fitnesse" are dropped, matching the target markup exactly.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$p1.InsertBefore("Acceptance Tests running on Node with Decaf and Coffeescript`r")
$d.Paragraphs(2).Range.Delete()

# ---------------------------------------------------------------------
# Paragraph 2: rewrite the intro blurb.
#   "on" -> "to"
#   "node.js" -> "Node.js"
#   "npm" -> "NPM"
#   "decaf" -> "Decaf"
#   + append the new sentence about Decaf/Coffeescript
# Rebuilt the same way, for the same w:proofErr-stripping reason (the
# original "npm" run was also spell-flagged).
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(2).Range
$apostrophe = [char]8217
$newIntro = "This session will serve as a quick intro to installing " `
    + "Node.js on windows and basic usage of NPM. Just enough needed " `
    + "to install and run Decaf. Decaf is the Slim port written in " `
    + "Coffeescript. We" + $apostrophe + "ll cover the basic setup and " `
    + "how to write Slim fixtures in Coffeescript. "
$p2.InsertBefore($newIntro + "`r")
$d.Paragraphs(3).Range.Delete()

# ---------------------------------------------------------------------
# The "_GoBack" bookmark originally sat at the very end of paragraph 2
# (Word drops it at the last edit point). In the target it has moved to
# just after "to" (between "to" and " installing"), reflecting that
# "on" -> "to" was the final edit made. Recreate it there.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$p2 = $d.Paragraphs(2).Range
$pos = $p2.Start + $p2.Text.IndexOf(" installing")
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)
